$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + new text value.
# Column D holds plain-number-looking strings (e.g. "231.43") that must
# stay text (matching the source inlineStr cells), so those are written
# with a temporary text NumberFormat and then restored to the default
# "Normal" style so no stray formatting is left behind.
$updates = @(
    @{ Cell = "D2"; Value = '34.901.67' },
    @{ Cell = "E2"; Value = '  -0.76%  ' },
    @{ Cell = "D3"; Value = '1.839.39' },
    @{ Cell = "E3"; Value = '  +1.29%  ' },
    @{ Cell = "D5"; Value = '231.43' },
    @{ Cell = "E5"; Value = '  -0.71%  ' },
    @{ Cell = "E6"; Value = '  +1.40%  ' },
    @{ Cell = "E7"; Value = '  -0.08%  ' },
    @{ Cell = "D8"; Value = '39.82' },
    @{ Cell = "E8"; Value = '  -2.84%  ' },
    @{ Cell = "D9"; Value = '0.329' },
    @{ Cell = "E9"; Value = '  +1.64%  ' },
    @{ Cell = "D10"; Value = '0.0686' },
    @{ Cell = "E10"; Value = '  +0.04%  ' },
    @{ Cell = "E11"; Value = '  -1.15%  ' },
    @{ Cell = "D12"; Value = '2.105.61' },
    @{ Cell = "E12"; Value = '  +1.33%  ' },
    @{ Cell = "D13"; Value = '11.47' },
    @{ Cell = "E13"; Value = '  +3.71%  ' },
    @{ Cell = "D14"; Value = '1.838.22' },
    @{ Cell = "E14"; Value = '  +1.39%  ' },
    @{ Cell = "D15"; Value = '0.673' },
    @{ Cell = "E15"; Value = '  +1.87%  ' },
    @{ Cell = "D16"; Value = '4.65' },
    @{ Cell = "E16"; Value = '  -0.10%  ' },
    @{ Cell = "D17"; Value = '34.890.30' },
    @{ Cell = "E17"; Value = '  -0.67%  ' },
    @{ Cell = "D18"; Value = '69.85' },
    @{ Cell = "E18"; Value = '  +0.33%  ' },
    @{ Cell = "D19"; Value = '0.0₃0787' },
    @{ Cell = "E19"; Value = '  -0.60%  ' },
    @{ Cell = "D20"; Value = '240.49' },
    @{ Cell = "E20"; Value = '  +0.58%  ' },
    @{ Cell = "E21"; Value = '  +2.27%  ' },
    @{ Cell = "D22"; Value = '4.70' },
    @{ Cell = "E22"; Value = '  +0.51%  ' },
    @{ Cell = "E23"; Value = '  -0.08%  ' },
    @{ Cell = "D24"; Value = '2.27' },
    @{ Cell = "E24"; Value = '  +0.15%  ' },
    @{ Cell = "D25"; Value = '171.13' },
    @{ Cell = "E25"; Value = '  -1.04%  ' },
    @{ Cell = "D26"; Value = '7.79' },
    @{ Cell = "E26"; Value = '  -0.85%  ' },
    @{ Cell = "E27"; Value = '  -0.48%  ' },
    @{ Cell = "E28"; Value = '  +2.51%  ' },
    @{ Cell = "E29"; Value = '  -4.85%  ' },
    @{ Cell = "E30"; Value = '  -0.14%  ' },
    @{ Cell = "E31"; Value = '  -0.11%  ' },
    @{ Cell = "D32"; Value = '3.95' },
    @{ Cell = "E32"; Value = '  -5.14%  ' },
    @{ Cell = "E33"; Value = '  -1.07%  ' },
    @{ Cell = "D34"; Value = '1.90' },
    @{ Cell = "E34"; Value = '  +7.25%  ' },
    @{ Cell = "E35"; Value = '  +6.88%  ' },
    @{ Cell = "D36"; Value = '1.42' },
    @{ Cell = "E36"; Value = '  +11.07%  ' },
    @{ Cell = "D37"; Value = '0.696' },
    @{ Cell = "E37"; Value = '  +2.28%  ' },
    @{ Cell = "E38"; Value = '  +6.26%  ' },
    @{ Cell = "D39"; Value = '91.13' },
    @{ Cell = "E39"; Value = '  -1.68%  ' },
    @{ Cell = "D40"; Value = '1.341.18' },
    @{ Cell = "E40"; Value = '  +2.32%  ' },
    @{ Cell = "E41"; Value = '  -0.15%  ' },
    @{ Cell = "D42"; Value = '14.81' },
    @{ Cell = "E42"; Value = '  +2.15%  ' },
    @{ Cell = "D43"; Value = '2.28' },
    @{ Cell = "E43"; Value = '  -0.94%  ' },
    @{ Cell = "E44"; Value = '  -2.95%  ' },
    @{ Cell = "E45"; Value = '  +0.04%  ' },
    @{ Cell = "D46"; Value = '6.29' },
    @{ Cell = "E46"; Value = '  -0.92%  ' },
    @{ Cell = "E47"; Value = '  +1.84%  ' },
    @{ Cell = "D48"; Value = '2.019.24' },
    @{ Cell = "E48"; Value = '  +1.36%  ' },
    @{ Cell = "E49"; Value = '  -0.08%  ' },
    @{ Cell = "B50"; Value = 'THORChain' },
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune' },
    @{ Cell = "D50"; Value = '3.38' },
    @{ Cell = "E50"; Value = '  +18.77%  ' },
    @{ Cell = "B51"; Value = 'Cronos' },
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = "D51"; Value = '0.0665' },
    @{ Cell = "E51"; Value = '  +1.93%  ' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell.Substring(0,1) -eq "D") {
        $range.NumberFormat = "@"
        $range.Value = $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
